# Update the "username" test-data value (B1, next to the "username" label in A1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "juli0o0dz1"

# Match the saved window geometry for this workbook window.
$win = $wb.Windows.Item(1)
$win.Left   = 40290
$win.Top    = 2265
$win.Width  = 18195
$win.Height = 13785

# Move the active selection on the sheet from C4 to C2.
$ws.Range("C2").Select()
